# Update "想去人数" (want-to-go count) values for a few events on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览": F2 1071 -> 1074, F4 1590 -> 1592, F6 37 -> 38
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1074
$wsExhibit.Range("F4").Value = 1592
$wsExhibit.Range("F6").Value = 38

# Sheet "全部类型": F2 1071 -> 1074, F4 1590 -> 1592, F7 37 -> 38
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1074
$wsAll.Range("F4").Value = 1592
$wsAll.Range("F7").Value = 38
